$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell G2 from "-" to "NA" (P6 response correction)
$ws.Range("G2").Value = "NA"

# Move active selection to E4
$ws.Range("E4").Select()
